$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.810.93"
$ws.Range("E2").Value = "  +2.49%  "

$ws.Range("D3").Value = "3.309.77"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.75"
$ws.Range("E5").Value = "  +3.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.16"
$ws.Range("E6").Value = "  -2.19%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +2.46%  "

$ws.Range("D9").Value = "3.306.07"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("E10").Value = "  +1.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.577"
$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.35"
$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("E13").Value = "  +2.96%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.843.26"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "628.41"
$ws.Range("E15").Value = "  +6.56%  "

$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").Value = "67.907.29"
$ws.Range("E17").Value = "  +2.70%  "

$ws.Range("E18").Value = "  +1.43%  "

$ws.Range("D19").Value = "3.317.73"
$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.90"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("E23").Value = "  -2.57%  "

$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.72"
$ws.Range("E25").Value = "  -1.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.98"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  +2.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("E28").Value = "  +1.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.24"
$ws.Range("E29").Value = "  +5.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.57"
$ws.Range("E30").Value = "  +1.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.70"
$ws.Range("E31").Value = "  +1.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "596.77"
$ws.Range("E32").Value = "  +6.33%  "

$ws.Range("D33").Value = "3.933.46"
$ws.Range("E33").Value = "  +3.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.94"
$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.104"
$ws.Range("E35").Value = "  +1.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.49"
$ws.Range("E36").Value = "  -5.29%  "

$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.83"
$ws.Range("E38").Value = "  -0.32%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.128"
$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.24"
$ws.Range("E40").Value = "  +3.36%  "

$ws.Range("E41").Value = "  +3.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "32.63"
$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.39"
$ws.Range("E43").Value = "  +0.40%  "

$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.337"
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("E46").Value = "  +0.43%  "

$ws.Range("E47").Value = "  +1.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("E49").Value = "  +12.66%  "

$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.32"
$ws.Range("E51").Value = "  +1.48%  "
